# Insert two new weekly records (rows 316-317) into the Cilantro price log.
# Excel pushes the existing rows 316:387 down to 318:389 and extends the
# sheet dimension automatically; we just need to fill in the values for the
# two freshly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 316 (shifts 316:387 -> 318:389).
$ws.Rows("316:317").Insert()

# Columns that are constant for every record in this block.
$constA = 9
$constB = "Vega Central Mapocho de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = 100112040
$constG = "Cilantro"
$constH = "Sin especificar"
$constI = "Primera"
$constO = "Regi$([char]0xF3)n Metropolitana"
$constR = "Hortaliza"

# New row 316: $/caja 36 atados record.
$ws.Cells.Item(316, 1).Value = $constA
$ws.Cells.Item(316, 2).Value = $constB
$ws.Cells.Item(316, 3).Value = $constC
$ws.Cells.Item(316, 4).Value = 44511
$ws.Cells.Item(316, 5).Value = $constE
$ws.Cells.Item(316, 6).Value = $constF
$ws.Cells.Item(316, 7).Value = $constG
$ws.Cells.Item(316, 8).Value = $constH
$ws.Cells.Item(316, 9).Value = $constI
$ws.Cells.Item(316, 10).Value = 43
$ws.Cells.Item(316, 11).Value = 5000
$ws.Cells.Item(316, 12).Value = 5000
$ws.Cells.Item(316, 13).Value = 5000
$ws.Cells.Item(316, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(316, 15).Value = $constO
$ws.Cells.Item(316, 16).Value = 139
$ws.Cells.Item(316, 17).Value = 36
$ws.Cells.Item(316, 18).Value = $constR

# New row 317: $/docena de atados record.
$ws.Cells.Item(317, 1).Value = $constA
$ws.Cells.Item(317, 2).Value = $constB
$ws.Cells.Item(317, 3).Value = $constC
$ws.Cells.Item(317, 4).Value = 44511
$ws.Cells.Item(317, 5).Value = $constE
$ws.Cells.Item(317, 6).Value = $constF
$ws.Cells.Item(317, 7).Value = $constG
$ws.Cells.Item(317, 8).Value = $constH
$ws.Cells.Item(317, 9).Value = $constI
$ws.Cells.Item(317, 10).Value = 160
$ws.Cells.Item(317, 11).Value = 10000
$ws.Cells.Item(317, 12).Value = 12000
$ws.Cells.Item(317, 13).Value = 11000
$ws.Cells.Item(317, 14).Value = "`$/docena de atados"
$ws.Cells.Item(317, 15).Value = $constO
$ws.Cells.Item(317, 16).Value = 3667
$ws.Cells.Item(317, 17).Value = 3
$ws.Cells.Item(317, 18).Value = $constR

# Make sure the date cells keep the date-formatted style used elsewhere in
# column D (copy the number format from the row above, which already has
# the correct style).
$ws.Range("D316").NumberFormat = $ws.Range("D315").NumberFormat
$ws.Range("D317").NumberFormat = $ws.Range("D315").NumberFormat
